# fix: fixed formatting when scrapping floating point numbers
#
# 1) A handful of "Razon social"/"Nombre Fantasia" entries that contained a
#    comma used as a separator between multiple people/names get the comma
#    turned into a period (and one of them additionally loses the dots in
#    the "S.H." abbreviation, becoming "SH").
# 2) Every amount in the "Importe" column (H2:H197) was stored as Spanish
#    locale formatted text (thousands separator "." and decimal separator
#    ",", e.g. "46.636,00"). They are normalized to plain decimal text
#    (e.g. "46636.00"): the thousands dots are removed and the decimal
#    comma becomes a dot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Name fields: replace the separator comma with a period ---------
$nameCells = @("E74", "E77", "E99", "E113", "E156", "F124")
foreach ($ref in $nameCells) {
    $ws.Range($ref).Replace(",", ".")
}

# Special case: also collapse "S.H." -> "SH" after the comma fix
$schab = $ws.Range("E164")
$schab.Replace(",", ".")
$schab.Replace("S.H.", "SH")

# --- 2) Importe column: reformat "1.234,56" -> "1234.56" ----------------
$importe = $ws.Range("H2:H197")

# Force text storage first so Excel doesn't reinterpret the reformatted
# value (which now looks like a plain number) as a numeric cell; this
# keeps it a shared string, exactly like the original data.
$importe.NumberFormat = "@"

$importe.Replace(".", "")
$importe.Replace(",", ".")

# Restore the default (unstyled) cell style now that the text is safely
# stored, so the cells end up just like the originals (no explicit style
# index attached to them).
$importe.Style = "Normal"
